# Update cryptos list with latest values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for cells whose new values look numeric, so Excel
# does not silently convert them to the Number type (matches the source
# workbook which stores every data cell as an inline/shared string).
$textRange = $ws.Range("D2:E51")
$textRange.NumberFormat = "@"

$ws.Range("D2").Value = '72.146.81'
$ws.Range("E2").Value = '  +4.09%  '
$ws.Range("D3").Value = '4.033.93'
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '521.20'
$ws.Range("E5").Value = '  -1.55%  '
$ws.Range("D6").Value = '147.09'
$ws.Range("E6").Value = '  +1.53%  '
$ws.Range("D7").Value = '0.716'
$ws.Range("E7").Value = '  +16.74%  '
$ws.Range("D8").Value = '4.024.35'
$ws.Range("E8").Value = '  +3.48%  '
$ws.Range("E9").Value = '  +0.15%  '
$ws.Range("D10").Value = '0.772'
$ws.Range("E10").Value = '  +7.16%  '
$ws.Range("E11").Value = '  +1.53%  '
$ws.Range("E12").Value = '  -2.03%  '
$ws.Range("D13").Value = '49.26'
$ws.Range("E13").Value = '  +17.13%  '
$ws.Range("D14").Value = '11.14'
$ws.Range("E14").Value = '  +8.54%  '
$ws.Range("D15").Value = '4.677.78'
$ws.Range("E15").Value = '  +3.45%  '
$ws.Range("D16").Value = '4.021.35'
$ws.Range("E16").Value = '  +3.36%  '
$ws.Range("E17").Value = '  +7.50%  '
$ws.Range("E18").Value = '  +1.62%  '
$ws.Range("E19").Value = '  +0.09%  '
$ws.Range("E20").Value = '  -0.51%  '
$ws.Range("D21").Value = '72.089.94'
$ws.Range("E21").Value = '  +4.08%  '
$ws.Range("D22").Value = '444.41'
$ws.Range("E22").Value = '  +4.71%  '
$ws.Range("D23").Value = '105.23'
$ws.Range("E23").Value = '  +19.74%  '
$ws.Range("D24").Value = '3.59'
$ws.Range("E24").Value = '  +5.94%  '
$ws.Range("D25").Value = '15.21'
$ws.Range("E25").Value = '  +7.22%  '
$ws.Range("D26").Value = '4.04'
$ws.Range("E26").Value = '  +0.81%  '
$ws.Range("D27").Value = '11.55'
$ws.Range("E27").Value = '  +0.91%  '
$ws.Range("D28").Value = '11.07'
$ws.Range("E28").Value = '  +4.27%  '
$ws.Range("D29").Value = '37.81'
$ws.Range("E29").Value = '  +3.89%  '
$ws.Range("E30").Value = '  +2.49%  '
$ws.Range("D31").Value = '3.28'
$ws.Range("E31").Value = '  +15.95%  '
$ws.Range("D32").Value = '13.77'
$ws.Range("E32").Value = '  +4.07%  '
$ws.Range("D33").Value = '0.131'
$ws.Range("E33").Value = '  +3.69%  '
$ws.Range("D34").Value = '675.72'
$ws.Range("E34").Value = '  -2.00%  '
$ws.Range("D35").Value = '6.78'
$ws.Range("E35").Value = '  +14.58%  '
$ws.Range("D36").Value = '67.76'
$ws.Range("E36").Value = '  -0.88%  '
$ws.Range("D37").Value = '42.50'
$ws.Range("E37").Value = '  +6.27%  '
$ws.Range("D38").Value = '0.0₃0860'
$ws.Range("E38").Value = '  -0.17%  '
$ws.Range("D39").Value = '0.427'
$ws.Range("E39").Value = '  -1.32%  '
$ws.Range("D40").Value = '3.52'
$ws.Range("E40").Value = '  +5.91%  '
$ws.Range("E41").Value = '  +1.51%  '
$ws.Range("E42").Value = '  +0.14%  '
$ws.Range("D43").Value = '0.0500'
$ws.Range("E43").Value = '  +3.48%  '
$ws.Range("D44").Value = '0.998'
$ws.Range("E44").Value = '  -0.25%  '
$ws.Range("D45").Value = '3.23'
$ws.Range("E45").Value = '  -0.09%  '
$ws.Range("D46").Value = '0.158'
$ws.Range("E46").Value = '  +12.37%  '
$ws.Range("E47").Value = '  -2.67%  '
$ws.Range("D48").Value = '3.50'
$ws.Range("E48").Value = '  +2.51%  '
$ws.Range("B49").Value = 'THORChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D49").Value = '9.58'
$ws.Range("E49").Value = '  +11.31%  '
$ws.Range("B50").Value = 'Stacks'
$ws.Range("C50").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D50").Value = '3.07'
$ws.Range("E50").Value = '  +2.90%  '
$ws.Range("E51").Value = '  +2.70%  '

# Remove the temporary text formatting so cell styles stay unchanged
$textRange.ClearFormats()

